{"js": "// Replace each two-digit multiplication expression in the grid with its\n// new value. Every source string is unique within the document, so a\n// plain case-sensitive search-and-replace on the whole body is safe.\nconst replacements = [\n  [\"69\u00d740=\", \"23\u00d746=\"],\n  [\"97\u00d760=\", \"29\u00d718=\"],\n  [\"72\u00d729=\", \"78\u00d771=\"],\n  [\"77\u00d791=\", \"17\u00d726=\"],\n  [\"71\u00d740=\", \"40\u00d752=\"],\n  [\"45\u00d760=\", \"69\u00d746=\"],\n  [\"18\u00d721=\", \"41\u00d751=\"],\n  [\"40\u00d795=\", \"90\u00d765=\"],\n  [\"33\u00d793=\", \"79\u00d767=\"],\n  [\"22\u00d739=\", \"54\u00d750=\"],\n  [\"55\u00d739=\", \"50\u00d718=\"],\n  [\"43\u00d717=\", \"46\u00d796=\"],\n  [\"67\u00d769=\", \"54\u00d724=\"],\n  [\"94\u00d785=\", \"98\u00d722=\"],\n  [\"98\u00d754=\", \"89\u00d716=\"],\n  [\"29\u00d761=\", \"13\u00d797=\"],\n  [\"40\u00d775=\", \"96\u00d752=\"],\n  [\"39\u00d763=\", \"70\u00d794=\"],\n  [\"48\u00d771=\", \"49\u00d781=\"],\n  [\"17\u00d763=\", \"86\u00d728=\"],\n  [\"93\u00d724=\", \"72\u00d715=\"],\n  [\"73\u00d787=\", \"59\u00d735=\"],\n  [\"60\u00d745=\", \"82\u00d794=\"],\n  [\"29\u00d755=\", \"21\u00d750=\"],\n  [\"34\u00d776=\", \"15\u00d712=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the grid with its\n# new value. Every source string is unique within the document (and none\n# of the new values collide with any of the old ones), so a single pass\n# over all paragraphs, matching each paragraph's leading text against the\n# map, is safe and order independent.\n$d = $word.ActiveDocument\n\n$map = @{\n  \"69\u00d740=\" = \"23\u00d746=\"\n  \"97\u00d760=\" = \"29\u00d718=\"\n  \"72\u00d729=\" = \"78\u00d771=\"\n  \"77\u00d791=\" = \"17\u00d726=\"\n  \"71\u00d740=\" = \"40\u00d752=\"\n  \"45\u00d760=\" = \"69\u00d746=\"\n  \"18\u00d721=\" = \"41\u00d751=\"\n  \"40\u00d795=\" = \"90\u00d765=\"\n  \"33\u00d793=\" = \"79\u00d767=\"\n  \"22\u00d739=\" = \"54\u00d750=\"\n  \"55\u00d739=\" = \"50\u00d718=\"\n  \"43\u00d717=\" = \"46\u00d796=\"\n  \"67\u00d769=\" = \"54\u00d724=\"\n  \"94\u00d785=\" = \"98\u00d722=\"\n  \"98\u00d754=\" = \"89\u00d716=\"\n  \"29\u00d761=\" = \"13\u00d797=\"\n  \"40\u00d775=\" = \"96\u00d752=\"\n  \"39\u00d763=\" = \"70\u00d794=\"\n  \"48\u00d771=\" = \"49\u00d781=\"\n  \"17\u00d763=\" = \"86\u00d728=\"\n  \"93\u00d724=\" = \"72\u00d715=\"\n  \"73\u00d787=\" = \"59\u00d735=\"\n  \"60\u00d745=\" = \"82\u00d794=\"\n  \"29\u00d755=\" = \"21\u00d750=\"\n  \"34\u00d776=\" = \"15\u00d712=\"\n}\n\nforeach ($p in $d.Paragraphs) {\n  $r = $p.Range\n  $t = $r.Text\n  foreach ($key in $map.Keys) {\n    if ($t.StartsWith($key)) {\n      $r.Text = $map[$key]\n      break\n    }\n  }\n}\n"}
